$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Ativacao date: 01/01/2012 -> 01/01/2023 (keep stored as text, like before) ---
$ws.Range("B8:C8").NumberFormat = "@"
$ws.Range("B8").Value = "01/01/2023"
$ws.Range("C8").Value = "01/01/2023"

# --- Objetivos: responsible professor changed ---
$ws.Range("B10").Value = "5840730 - Antonio Jefferson da Silva Machado"
$ws.Range("C10").Value = "5840730 - Antonio Jefferson da Silva Machado"

# --- Objectives (English) text newly added in row 11 (B/C cells did not exist before) ---
$ws.Range("B11").WrapText = $true
$ws.Range("B11").VerticalAlignment = -4160
$ws.Range("B11").Font.Bold = $false
$ws.Range("B11").Value = "Provide the student with an overview of the area of Solid State Physics, with emphasis on fundamental ideas and general concepts, such as electron gas, elementary excitations, band structure, etc. The course should be rich in experimental results that illustrate general principles and behaviors of solids (eg, behavior of physical quantities with temperature)."

$ws.Range("C11").WrapText = $true
$ws.Range("C11").VerticalAlignment = -4160
$ws.Range("C11").Font.Color = 255
$ws.Range("C11").Value = "Provide the student with an overview of the area of Solid State Physics, with emphasis on fundamental ideas and general concepts, such as electron gas, elementary excitations, band structure, etc. The course should be rich in experimental results that illustrate general principles and behaviors of solids (eg, behavior of physical quantities with temperature)."

# --- Programa resumido: responsible professor (same as row 10) ---
$ws.Range("B13").Value = "5840730 - Antonio Jefferson da Silva Machado"
$ws.Range("C13").Value = "5840730 - Antonio Jefferson da Silva Machado"

# --- Short syllabus (English) text newly added in row 14 ---
$ws.Range("B14").WrapText = $true
$ws.Range("B14").VerticalAlignment = -4160
$ws.Range("B14").Font.Bold = $false
$ws.Range("B14").Value = "Crystal structure and bonds. Lattice vibrations, phonons and thermal properties. Free electron Fermi gas. Power bands. Semiconductors. Fermi metals and surfaces."

$ws.Range("C14").WrapText = $true
$ws.Range("C14").VerticalAlignment = -4160
$ws.Range("C14").Font.Color = 255
$ws.Range("C14").Value = "Crystal structure and bonds. Lattice vibrations, phonons and thermal properties. Free electron Fermi gas. Power bands. Semiconductors. Fermi metals and surfaces."

# --- Programa: professor changed ---
$ws.Range("B15").Value = "5840726 - Cristina Bormio Nunes"
$ws.Range("C15").Value = "5840726 - Cristina Bormio Nunes"

# --- Syllabus (English) text newly added in row 16 ---
$ws.Range("B16").WrapText = $true
$ws.Range("B16").VerticalAlignment = -4160
$ws.Range("B16").Font.Bold = $false
$ws.Range("B16").Value = "¨ Structure of crystals.¨ Crystal diffraction and the reciprocal lattice.¨ Bonds in crystals: ionic crystals and covalent crystals¨ Elastic constants and elastic waves.¨ Crystal vibrations. phonons¨ Fermi gas: free electron model; movement in magnetic fields.¨ Energy bands. Bloch functions.¨ Semiconductor crystals."

$ws.Range("C16").WrapText = $true
$ws.Range("C16").VerticalAlignment = -4160
$ws.Range("C16").Font.Color = 255
$ws.Range("C16").Value = "¨ Structure of crystals.¨ Crystal diffraction and the reciprocal lattice.¨ Bonds in crystals: ionic crystals and covalent crystals¨ Elastic constants and elastic waves.¨ Crystal vibrations. phonons¨ Fermi gas: free electron model; movement in magnetic fields.¨ Energy bands. Bloch functions.¨ Semiconductor crystals."

# --- Norma de recuperacao: grading policy text simplified ---
$ws.Range("B20").Value = "Média aritmética de duas provas com mesmo peso."
$ws.Range("C20").Value = "Média aritmética de duas provas com mesmo peso."
